$d = $word.ActiveDocument

function Replace-ExactText($oldText, $newText) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Text = $oldText
    $rng.Find.Forward = $true
    $rng.Find.Wrap = 1
    $rng.Find.MatchCase = $true
    $rng.Find.MatchWholeWord = $false
    $rng.Find.MatchWildcards = $false
    $found = $rng.Find.Execute()
    if (-not $found) {
        throw "Text not found: $oldText"
    }
    $rng.Text = $newText
}

# Title
Replace-ExactText "Quantum Entanglement: Unveiling the Mysteries of Interconnectedness" "Biology: Life's Evolving Tapestry"

# Author name
Replace-ExactText "Eleanor Bergstrom" "Hazel Scott"

# Email (username part and domain suffix part, separated by a "." run)
Replace-ExactText "ebergstrom@berkeley" "hscott@humbleton"
Replace-ExactText "edu" "academy"

# Body paragraph, first visual line
Replace-ExactText "In the realm of quantum mechanics, the universe reveals a profound enigma known as quantum entanglement. This extraordinary phenomenon challenges our classical understanding of reality, introducing the notion of interconnectedness between particles separated by vast distances. The intricate dance of entangled particles transcends time and space, defying traditional notions of causality and ushering us into a realm of astonishing possibilities. As we delve into the depths of quantum entanglement, we embark on a captivating journey to unravel the mysteries that lie at the heart of our physical universe." "Biology, an intricate symphony of life, unveils the enchanting saga of living organisms that inhabit our planet. It is a realm where countless species engage in a delicate dance of interactions, defining ecosystems and shaping Earth's history. From the smallest microscopic entity to the vast expanse of rainforests, biology orchestrates an incredible narrative of adaptation and resilience."

# Body paragraph, second visual line
Replace-ExactText "In the tapestry of quantum mechanics, particles can exhibit remarkable correlations that defy explanation based on classical physics. Entangled particles, once brought into contact, become inextricably linked, regardless of the distance that may subsequently separate them. The actions performed on one entangled particle instantaneously affect the state of its distant counterpart. This profound phenomenon, defying the constraints of locality, has captivated the imaginations of scientists and philosophers alike, challenging our fundamental understanding of reality." "Embarking on this journey of discovery, we begin by understanding the basic building blocks that constitute all living organisms - cells. Within these microscopic worlds, genetic information encoded in DNA guides the symphony of cellular processes, ensuring life's continuity across generations. We delve into the intricacies of photosynthesis, the magical process that converts sunlight into energy, fueling plant growth and sustaining the food chain's very foundation. Our explorations extend to the incredible diversity of life forms, from single-celled organisms like bacteria to complex, multicellular creatures like humans. We unravel the intricate tapestry of ecosystems, where intricate webs of relationships weave together organisms and their environments."

# Body paragraph, third visual line
Replace-ExactText "Intriguingly, the interconnectedness of entangled particles transcends the limitations of space and time. Measurements performed on one particle instantaneously impact the properties of its entangled partner, even if they are separated by vast cosmic distances. This nonlocal connection, known as quantum nonlocality, has profound implications for our understanding of the universe. It suggests the existence of a deeper level of reality, beyond the realm of our everyday experiences, where particles communicate and influence each other instantaneously." "Further, we explore the amazing adaptations that enable organisms to thrive in various environments, showcasing nature's boundless creativity. From the stunning camouflage techniques that help animals blend with their surroundings to the sophisticated echolocation abilities of bats, we unravel the secrets of survival in ecosystems that challenge our imaginations. We examine the concept of homeostasis, the body's ability to maintain internal balance in the face of external changes. We also delve into genetics, deciphering the intricate code of life etched within DNA, unlocking the secrets of genetic inheritance and variation."

# Summary paragraph
Replace-ExactText "Quantum entanglement, a perplexing phenomenon in the realm of quantum mechanics, intertwines the fates of particles across vast distances, blurring the lines between time and space. This phenomenon, defying intuition and challenging classical notions of causality, has profound implications for our understanding of the universe, suggesting interconnectedness at the deepest levels of reality. The instantaneous transfer of information between entangled particles defies locality, hinting at a nonlocal connection that transcends the constraints of space and time. As we delve further into the mysteries of quantum entanglement, we may uncover hidden truths about the fundamental nature of our physical universe, inviting us to reconceptualize reality itself." "In the grand theater of life, biology plays the starring role, orchestrating the delicate interplay between organisms and their environments. Through its enchanting stories of adaptation, resilience, and diversity, biology unfolds before us a rich tapestry of knowledge, captivating and inspiring young minds to explore the wonders of life on Earth."

# Add an empty paragraph at the very end of the document body (after the Summary paragraph)
$endRng = $d.Content
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()
